$p = $ppt.ActivePresentation

# --- 1) Update the "Date Placeholder" field text on the slide master and
#        every slide layout (datetimeFigureOut field: 09/08/2024 -> 12/11/2025).
function Set-DateFieldText($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

$master = $p.SlideMaster
Set-DateFieldText $master.Shapes "12/11/2025"
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Set-DateFieldText $layout.Shapes "12/11/2025"
}

# --- 2) Slide 1: move the textbox right and change its text.
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(1)

# off x="-396665" -> x="-196640" (EMU); 1 pt = 12700 EMU.
$shp.Left = -15.483466

# Text content: 中文 -> 简中
$shp.TextFrame.TextRange.Text = "简中"

# Re-setting the text re-triggers the shape's auto-fit, nudging the height
# by a few EMU; restore the original extent (cy="4696222") so only the
# intended attributes (off/text) change.
$shp.Height = 369.78128
